# Session plan wizard field-name fix.
# Applies the content updates described by the commit:
# "Fix: Update document generator to match wizard form field names -
#  Session plans and schemes now use correct data from forms"
#
# Note: in the Find/Replace text strings below, "^l" is Word's own
# wildcard notation for a manual line break (<w:br/>) - it must be the
# literal two characters "^" + "l", NOT an embedded vertical-tab char,
# since Find.Execute's replacement parser only recognises the caret
# escape form.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Header / meta fields ---------------------------------------------
Replace-Text "Sector :    SGGS" "Sector :    ICT & MULTIMEDIA"
Replace-Text "Sub-sector: HH" "Sub-sector: Software Development"
Replace-Text "TERM : I" "TERM : Term 1"
Replace-Text "Module(Code&Name): TEST301 Test Module" "Module(Code&Name): SWDPR301: Analyze project requirements"
Replace-Text "Week : I" "Week : Week 1"
Replace-Text "Class(es): 1" "Class(es): L3SD-A"

# "Lead Trainer's name : SJJS" -> "... TUYISINGIZE Leonard" contains an
# apostrophe in unchanged text ("Trainer's"). Find/Replace re-types the
# whole match and Word's AutoCorrect turns straight quotes into curly
# ones in freshly (re)typed text, so set the paragraph Range.Text
# directly instead - that preserves the literal apostrophe.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Lead Trainer*SJJS*") {
        $p.Range.Text = "Lead Trainer's name : TUYISINGIZE Leonard"
    }
}

# --- Learning outcome / indicative contents ----------------------------
Replace-Text "1. Test learning outcome" "1. Identify customer needs^l2. Apply data gathering techniques"
Replace-Text "1.1 Test content" "1.1 Data gathering^l1.2 Communication process^l1.3 Customer pain points"

# --- Topic / range / duration ------------------------------------------
Replace-Text "Topic of the session: Push the latest code to github pages" "Topic of the session: Identification of requirements gathering methodologies"
Replace-Text "Key concepts^lTest methods" "Level 3"
Replace-Text "Duration of the session: 55min" "Duration of the session: 80min"

# --- Objectives ----------------------------------------------------------
Replace-Text "Define key concepts^lApply test methods" "1. Identify customer needs^l2. Apply data gathering techniques"

# --- Facilitation technique ----------------------------------------------
Replace-Text "Facilitation technique(s):   JIGSAW" "Facilitation technique(s):   Jigsaw"

# --- Introduction step -----------------------------------------------------
Replace-Text "Greets and makes roll calls" "Greets and makes roll calls^lIntroduces topic and objectives"
Replace-Text "Attendance sheet^lPPT^lProjector" "Attendance sheet^lPPT^lProjector^lWhiteboard"

# --- Step 1 ---------------------------------------------------------------
Replace-Text "Step 1: Introduction to concepts" "Step 1: Introduction to topic"
Replace-Text "Explains key concepts" "Explains key concepts and demonstrates examples"

# --- Step 2 ---------------------------------------------------------------
Replace-Text "Demonstrates practical examples" "Guides learners through hands-on activities"

# Resources for Step 1, Step 2 and Step 3 all share the identical
# "Computer/Projector/PPT" text and all three gain "Learning materials" -
# a single ReplaceAll covers all three occurrences.
Replace-Text "Computer^lProjector^lPPT" "Computer^lProjector^lPPT^lLearning materials"

# Step 1 & Step 2 durations ("20 minutes" x2) and Step 3 duration
# ("10 minutes") all become "16" + line break + "minutes".
Replace-Text "20 minutes" "16^lminutes"
Replace-Text "10 minutes" "16^lminutes"

# --- Step 3 ---------------------------------------------------------------
Replace-Text "Step 3: Group work" "Step 3: Group work and practice"
Replace-Text "Assigns group tasks" "Assigns tasks and monitors learner progress"

# --- Conclusion: summary ---------------------------------------------------
Replace-Text "Trainer involves learners to summarize" "Trainer involves learners to summarize key points learned"

# --- Assessment -------------------------------------------------------------
Replace-Text "Gives assessment questions" "Provides assessment questions or assignment"
Replace-Text "Assessment sheets" "Assessment sheets^lQuestion papers"

# --- Evaluation ---------------------------------------------------------------
Replace-Text "Involves learners in evaluation" "Collects feedback from learners"
Replace-Text "Self-assessment form" "Self-assessment form^lFeedback forms"

# --- References / appendices ---------------------------------------------------
Replace-Text "RTB Curriculum Guide 2024" "RTB Curriculum Guidelines^lModule Learning Materials"
Replace-Text "Appendices: PPT, Task Sheets, Assessment" "Appendices: PPT, Task Sheets, Assessment, Learning Materials"
